$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 71430350
$ws.Range("I40").Value = 2372.25
$ws.Range("J40").Value = 166667660
$ws.Range("K40").Value = 2372.25
$ws.Range("L40").Value = 166667660
$ws.Range("M40").Value = -2197.25
$ws.Range("N40").Value = -166668010
$ws.Range("H51").Value = 8185.2354
$ws.Range("I51").Value = 12860
$ws.Range("J51").Value = 6237.4165
$ws.Range("K51").Value = 12860
$ws.Range("L51").Value = 6237.4165
$ws.Range("M51").Value = -12376
$ws.Range("N51").Value = -7205.4165
$ws.Range("H97").Value = 6387.6
$ws.Range("J97").Value = 6387.6
$ws.Range("L97").Value = 19162.8
$ws.Range("N97").Value = -20154.8
$ws.Range("H101").Value = 630.9286
$ws.Range("I101").Value = 1984
$ws.Range("J101").Value = 405.41666
$ws.Range("K101").Value = 5952
$ws.Range("L101").Value = 1216.24998
$ws.Range("M101").Value = -4330
$ws.Range("N101").Value = -4460.249980000001
$ws.Range("H112").Value = 3717.0625
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 3717.0625
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 11151.1875
$ws.Range("M112").Value = ""
$ws.Range("N112").Value = -13367.1875
$ws.Range("H135").Value = 4080.3635
$ws.Range("I135").Value = 610.625
$ws.Range("J135").Value = 13333
$ws.Range("K135").Value = 5495.625
$ws.Range("L135").Value = 119997
$ws.Range("M135").Value = -2960.625
$ws.Range("N135").Value = -125067
$ws.Range("H138").Value = 3910
$ws.Range("I138").Value = 2994.3845
$ws.Range("J138").Value = 4231.7026
$ws.Range("K138").Value = 8983.1535
$ws.Range("L138").Value = 12695.1078
$ws.Range("M138").Value = -3843.1535
$ws.Range("N138").Value = -22975.1078
$ws.Range("H141").Value = 3553.5454
$ws.Range("I141").Value = 2881.878
$ws.Range("K141").Value = 8645.634
$ws.Range("M141").Value = -3465.634

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 62765.25
$ws.Range("J24").Value = 62765.25
$ws.Range("L24").Value = 62765.25
$ws.Range("N24").Value = -63513.25
$ws.Range("H53").Value = 24949.5
$ws.Range("I53").Value = 9900
$ws.Range("J53").Value = 39999
$ws.Range("K53").Value = 9900
$ws.Range("L53").Value = 39999
$ws.Range("M53").Value = -9218
$ws.Range("N53").Value = -41363
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").Value = ""
$ws.Range("H60").Value = 89990.75
$ws.Range("I60").Value = 89990.75
$ws.Range("K60").Value = 89990.75
$ws.Range("M60").Value = -89257.75
$ws.Range("H88").Value = 2991.3333
$ws.Range("I88").Value = 2225
$ws.Range("J88").Value = 3374.5
$ws.Range("K88").Value = 2225
$ws.Range("L88").Value = 3374.5
$ws.Range("M88").Value = -1819
$ws.Range("N88").Value = -4186.5
$ws.Range("H91").Value = 2991.3333
$ws.Range("I91").Value = 2225
$ws.Range("J91").Value = 3374.5
$ws.Range("K91").Value = 2225
$ws.Range("L91").Value = 3374.5
$ws.Range("M91").Value = -821
$ws.Range("N91").Value = -6182.5
$ws.Range("H96").Value = 34147.6
$ws.Range("J96").Value = 34147.6
$ws.Range("L96").Value = 34147.6
$ws.Range("N96").Value = -39639.6
$ws.Range("H100").Value = 62765.25
$ws.Range("J100").Value = 62765.25
$ws.Range("L100").Value = 62765.25
$ws.Range("N100").Value = -64929.25
$ws.Range("H132").Value = 2071.0571
$ws.Range("I132").Value = 2186.6775
$ws.Range("K132").Value = 6560.032499999999
$ws.Range("M132").Value = -4030.032499999999

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3499.0334
$ws.Range("I86").Value = 2499.75
$ws.Range("J86").Value = 5497.6
$ws.Range("K86").Value = 2499.75
$ws.Range("L86").Value = 5497.6
$ws.Range("M86").Value = -1376.75
$ws.Range("N86").Value = -7743.6
$ws.Range("H89").Value = 3499.0334
$ws.Range("I89").Value = 2499.75
$ws.Range("J89").Value = 5497.6
$ws.Range("K89").Value = 12498.75
$ws.Range("L89").Value = 27488
$ws.Range("M89").Value = -6882.75
$ws.Range("N89").Value = -38720

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 5000
$ws.Range("I2").Value = 5000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 5000
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -4887
$ws.Range("N2").Value = ""
$ws.Range("H33").Value = 11531
$ws.Range("I33").Value = 11531
$ws.Range("K33").Value = 11531
$ws.Range("M33").Value = -11152
$ws.Range("H36").Value = 60048
$ws.Range("I36").Value = 60048
$ws.Range("K36").Value = 60048
$ws.Range("M36").Value = -59660
$ws.Range("H40").Value = 60048
$ws.Range("I40").Value = 60048
$ws.Range("K40").Value = 60048
$ws.Range("M40").Value = -59888
$ws.Range("H132").Value = 1695.625
$ws.Range("I132").Value = 1581.4849
$ws.Range("J132").Value = 2233.7144
$ws.Range("K132").Value = 4744.4547
$ws.Range("L132").Value = 6701.1432
$ws.Range("M132").Value = -2214.4547
$ws.Range("N132").Value = -11761.1432

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1793.4667
$ws.Range("I12").Value = 1954
$ws.Range("J12").Value = 1686.4445
$ws.Range("K12").Value = 5862
$ws.Range("L12").Value = 5059.333500000001
$ws.Range("M12").Value = -5689
$ws.Range("N12").Value = -5405.333500000001
$ws.Range("H132").Value = 2732
$ws.Range("J132").Value = 3779.2222
$ws.Range("L132").Value = 34012.99980000001
$ws.Range("N132").Value = -39072.99980000001
$ws.Range("H134").Value = 4367.615
$ws.Range("I134").Value = 1953.8334
$ws.Range("J134").Value = 33333
$ws.Range("K134").Value = 5861.5002
$ws.Range("L134").Value = 99999
$ws.Range("M134").Value = -791.5002000000004
$ws.Range("N134").Value = -110139
$ws.Range("H136").Value = 11015
$ws.Range("I136").Value = 3575.6667
$ws.Range("K136").Value = 10727.0001
$ws.Range("M136").Value = -5627.000100000001

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 209.27272
$ws.Range("I2").Value = 87.59999999999999
$ws.Range("J2").Value = 310.66666
$ws.Range("K2").Value = 87.59999999999999
$ws.Range("L2").Value = 310.66666
$ws.Range("M2").Value = 25.40000000000001
$ws.Range("N2").Value = -536.66666
$ws.Range("H53").Value = 40019
$ws.Range("I53").Value = 30039
$ws.Range("K53").Value = 30039
$ws.Range("M53").Value = -29408
$ws.Range("H55").Value = 52534.5
$ws.Range("I55").Value = 20069
$ws.Range("J55").Value = 85000
$ws.Range("K55").Value = 20069
$ws.Range("L55").Value = 85000
$ws.Range("M55").Value = -19742
$ws.Range("N55").Value = -85654
$ws.Range("H59").Value = 61166
$ws.Range("I59").Value = 48999
$ws.Range("K59").Value = 48999
$ws.Range("M59").Value = -48416
$ws.Range("H132").Value = 2334467
$ws.Range("I132").Value = 3815.6875
$ws.Range("J132").Value = 12988873
$ws.Range("K132").Value = 11447.0625
$ws.Range("L132").Value = 38966619
$ws.Range("M132").Value = -8917.0625
$ws.Range("N132").Value = -38971679

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3125.3076
$ws.Range("I22").Value = 3458
$ws.Range("J22").Value = 2593
$ws.Range("K22").Value = 3458
$ws.Range("L22").Value = 2593
$ws.Range("M22").Value = -3163
$ws.Range("N22").Value = -3183
$ws.Range("H27").Value = 3125.3076
$ws.Range("I27").Value = 3458
$ws.Range("J27").Value = 2593
$ws.Range("K27").Value = 3458
$ws.Range("L27").Value = 2593
$ws.Range("M27").Value = -3351
$ws.Range("N27").Value = -2807
$ws.Range("H68").Value = 1489717.6
$ws.Range("I68").Value = 2084620.2
$ws.Range("J68").Value = 2461.25
$ws.Range("K68").Value = 2084620.2
$ws.Range("L68").Value = 2461.25
$ws.Range("M68").Value = -2083871.2
$ws.Range("N68").Value = -3959.25
$ws.Range("H71").Value = 1489717.6
$ws.Range("I71").Value = 2084620.2
$ws.Range("J71").Value = 2461.25
$ws.Range("K71").Value = 10423101
$ws.Range("L71").Value = 12306.25
$ws.Range("M71").Value = -10419357
$ws.Range("N71").Value = -19794.25
$ws.Range("H132").Value = 2467.3088
$ws.Range("I132").Value = 1594.7407
$ws.Range("K132").Value = 4784.2221
$ws.Range("M132").Value = -2254.2221
$ws.Range("H136").Value = 3961.037
$ws.Range("I136").Value = 2472.7273
$ws.Range("K136").Value = 7418.1819
$ws.Range("M136").Value = -4868.1819
$ws.Range("H140").Value = 289429
$ws.Range("J140").Value = 289429
$ws.Range("L140").Value = 289429
$ws.Range("N140").Value = -299789

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").Value = ""
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").Value = ""
$ws.Range("H62").Value = 1225113.2
$ws.Range("J62").Value = 2287338.5
$ws.Range("L62").Value = 2287338.5
$ws.Range("N62").Value = -2288586.5
$ws.Range("H65").Value = 1225113.2
$ws.Range("J65").Value = 2287338.5
$ws.Range("L65").Value = 11436692.5
$ws.Range("N65").Value = -11442932.5
$ws.Range("H107").Value = 3813.55
$ws.Range("I107").Value = 2254.6956
$ws.Range("J107").Value = 5922.5884
$ws.Range("K107").Value = 6764.0868
$ws.Range("L107").Value = 17767.7652
$ws.Range("M107").Value = -4844.0868
$ws.Range("N107").Value = -21607.7652
$ws.Range("H132").Value = 1255.582
$ws.Range("I132").Value = 1124.6666
$ws.Range("J132").Value = 1799.3846
$ws.Range("K132").Value = 3373.9998
$ws.Range("L132").Value = 5398.1538
$ws.Range("M132").Value = -843.9998000000001
$ws.Range("N132").Value = -10458.1538
